$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "277.45"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.75%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.49"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.30%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.846"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.03%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06371"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.06%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.000"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.11%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.365"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "6.36%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8755"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.59%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1516"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.32%"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.93%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07521"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.68%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02933"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.94%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08963"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.81%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001566"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.47%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006389"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.63%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006082"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.18%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.477"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.74%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.301"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.09%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.06%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3147"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.92%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.55%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.909"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.35%"
$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.1519"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "10.04%"
$ws.Range("B24").Value = "CoinExToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04407"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.38%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001176"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.71%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "1.74%"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-1.40%"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "-0.41%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04069"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.16%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006820"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-3.20%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1414"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "20.75%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-6.25%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01167"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "0.18%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "3.31%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "4.79%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-19.42%"
